# Crypto price list refresh (Tue Oct 31 15:22:35 UTC 2023, GitHub Actions).
# Updates the "Price" (D) and "Volume(1h)" (E) columns for every coin row,
# and rotates rows 40-42 (Aave / HuobiToken / MXToken) to reflect the
# source feed's new ranking order, together with their refreshed figures.
#
# D-column values that look like plain numbers ("224.98", "11.07", ...) are
# written with a leading "'" so Excel keeps them as text (matching the
# workbook's existing text-typed Price column) instead of silently
# re-typing the cell as a Number; ".Style" is then reset to strip the
# quote-prefix formatting Excel would otherwise stamp on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.303.94"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.793.11"
$ws.Range("E3").Value = "  -1.56%  "
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "'224.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.52%  "
$ws.Range("D6").Value = "'0.588"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.58%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'35.83"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.44%  "
$ws.Range("D9").Value = "'0.289"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.43%  "
$ws.Range("D10").Value = "'0.0670"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.33%  "
$ws.Range("D11").Value = "'0.0960"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.76%  "
$ws.Range("D12").Value = "2.059.35"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "'11.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.80%  "
$ws.Range("D14").Value = "1.806.42"
$ws.Range("E14").Value = "  -0.97%  "
$ws.Range("D15").Value = "'0.626"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.28%  "
$ws.Range("D16").Value = "34.307.05"
$ws.Range("E16").Value = "  -1.15%  "
$ws.Range("D17").Value = "'4.33"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.54%  "
$ws.Range("D18").Value = "'68.62"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.95%  "
$ws.Range("D19").Value = "'239.09"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "0.0₃0764"
$ws.Range("E20").Value = "  -4.68%  "
$ws.Range("E21").Value = "  +0.12%  "
$ws.Range("D22").Value = "'11.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.90%  "
$ws.Range("D23").Value = "'4.07"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.80%  "
$ws.Range("D24").Value = "'2.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.07%  "
$ws.Range("D25").Value = "'169.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.27%  "
$ws.Range("D26").Value = "'7.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.62%  "
$ws.Range("D27").Value = "'17.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$ws.Range("D28").Value = "'0.119"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +0.24%  "
$ws.Range("D30").Value = "'1.23"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.38%  "
$ws.Range("D31").Value = "'3.75"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.90%  "
$ws.Range("D32").Value = "'3.85"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.78%  "
$ws.Range("D33").Value = "'0.0507"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.57%  "
$ws.Range("D34").Value = "'1.75"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.09%  "
$ws.Range("D35").Value = "1.355.15"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "'0.638"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -6.05%  "
$ws.Range("D37").Value = "'1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("D38").Value = "'2.32"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -10.02%  "
$ws.Range("D39").Value = "'0.0183"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.46%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'2.42"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.88%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "'2.78"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'80.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.00%  "
$ws.Range("D43").Value = "'0.922"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.28%  "
$ws.Range("E44").Value = "  +5.32%  "
$ws.Range("D45").Value = "'13.07"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.09%  "
$ws.Range("D46").Value = "'0.0496"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.54%  "
$ws.Range("D47").Value = "1.958.31"
$ws.Range("E47").Value = "  -1.35%  "
$ws.Range("D48").Value = "'5.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.39%  "
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").Value = "'101.35"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.83%  "
$ws.Range("D51").Value = "0.0₆0123"
$ws.Range("E51").Value = "  -6.09%  "
